$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the pattern-element mapping cells to reflect the renamed/merged
# function & class names (cleanup of naming in the lab code).
$ws.Range("E12").Value2 = "Perspective, ImagePerspectivePackage"
$ws.Range("E15").Value2 = "ZoomInCommand, ZoomOutCommand,  TranslateFreeCommand, LoadCommand, SerializeCommand, DeserializeCommand"
$ws.Range("E13").Value2 = "ImagePanel, MainPanel"

# The comment/description row shrinks now that the text is shorter.
$ws.Rows.Item(15).RowHeight = 46.5

# Update the view: scroll position and current selection.
$ws.Range("D7:E13").Select()
